$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Value" (column C) figures for rows 2018-2497 (GESS model refresh)
$newValues = @(
    6199.8,
    6198.6,
    6229.9,
    6215.8,
    5914.3,
    5888.9,
    5865.2,
    5882.9,
    5880.4,
    5816,
    5813.5,
    5798.8,
    5803,
    5789.7,
    5781.9,
    5742.5,
    5803.1,
    5795.5,
    5796.8,
    5812.1,
    5789.3,
    5790.5,
    5800.8,
    5713.6,
    5747.5,
    5750.2,
    5779.5,
    5767,
    6287.3,
    6279.4,
    6420.2,
    6457.8,
    6396,
    6487.5,
    6596.5,
    6663.8,
    6479,
    6508.4,
    6553,
    6603,
    6507,
    6530.9,
    6566.7,
    6607.3,
    6173.8,
    6195.5,
    6253.5,
    6281.7,
    6371.9,
    6393.4,
    6381.2,
    6408.9,
    6210.9,
    6188.2,
    6185.8,
    6200.3,
    6211.4,
    6181.8,
    6151.8,
    6131.1,
    6093.5,
    6056.4,
    6029.3,
    5975,
    6412.7,
    6482.3,
    6462.5,
    6444.8,
    6948.9,
    7001.8,
    7050.2,
    7044,
    7211.5,
    7199.7,
    7226.1,
    7251.9,
    7221.2,
    7236.5,
    7232.7,
    7233.6,
    7044.5,
    7033.3,
    6968.3,
    6991.8,
    6902.8,
    6883,
    6830.2,
    6798.4,
    6591.3,
    6534,
    6486.2,
    6473.6,
    6364.6,
    6307.1,
    6239.6,
    6191.7,
    6004.1,
    5994.4,
    5983.9,
    5952.1,
    5589,
    5559,
    5519,
    5509,
    5399,
    5389,
    5379,
    5379,
    5361,
    5361,
    5371,
    5371,
    5375,
    5395,
    5425,
    5485,
    5626,
    5736,
    5856,
    6006,
    6239,
    6429,
    6619,
    6819,
    7096,
    7286,
    7466,
    7636,
    7975,
    8105,
    8215,
    8305,
    8381,
    8431,
    8461,
    8461,
    8386,
    8346,
    8306,
    8246,
    8164,
    8104,
    8044,
    7994,
    7945,
    7915,
    7895,
    7875,
    7874,
    7874,
    7864,
    7864,
    7864,
    7854,
    7844,
    7844,
    7834,
    7834,
    7844,
    7864,
    7746,
    7796,
    7856,
    7926,
    8070,
    8140,
    8200,
    8260,
    8237,
    8277,
    8307,
    8337,
    8372,
    8362,
    8332,
    8292,
    8307,
    8247,
    8157,
    8057,
    7880,
    7750,
    7630,
    7510,
    7440,
    7300,
    7140,
    7000,
    6954,
    6814,
    6694,
    6584,
    6272,
    6202,
    6152,
    6112,
    6096,
    6071,
    6035.8,
    5996.7,
    5965.2,
    5932.5,
    5937.3,
    5883.7,
    5934.3,
    5927.8,
    5912.4,
    5907.9,
    5886.8,
    5893.9,
    5884.1,
    5907.9,
    5836.6,
    5879.9,
    5904.8,
    5912.9,
    5917.8,
    6001.2,
    6098.7,
    6161,
    6918.9,
    6970.9,
    7099.7,
    7157.6,
    7105.7,
    7175.4,
    7204.1,
    7297.3,
    7242.8,
    7248.8,
    7285.1,
    7298,
    7274.4,
    7278.5,
    7320.1,
    7347.8,
    6990.9,
    7025.6,
    7054.3,
    7115.4,
    7113,
    7168,
    7214.6,
    7279.2,
    7251.6,
    7253.9,
    7306.1,
    7348.9,
    7441.9,
    7443.3,
    7435.1,
    7442.8,
    7595,
    7530.6,
    7484.4,
    7425.6,
    7893.1,
    7891.6,
    7896,
    7895.5,
    7870,
    7870.5,
    7893,
    7893.9,
    8229.799999999999,
    8211.1,
    8245.700000000001,
    8288.6,
    8373.6,
    8440.299999999999,
    8496,
    8534.9,
    8263.9,
    8298.1,
    8350.299999999999,
    8395.299999999999,
    8274.799999999999,
    8226.700000000001,
    8225.200000000001,
    8224.1,
    8205.5,
    8181.7,
    8111,
    8108.7,
    8089.8,
    8046.7,
    8008.9,
    7985,
    6182,
    6132,
    6082,
    6032,
    6026,
    6001,
    5945.8,
    5906.7,
    5895.2,
    5852.5,
    5877.3,
    5843.7,
    5914.3,
    5917.8,
    5922.4,
    5927.9,
    5896.8,
    5903.9,
    5894.1,
    5897.9,
    5826.6,
    5849.9,
    5874.8,
    5882.9,
    5907.8,
    6001.2,
    6108.7,
    6181,
    6878.9,
    6920.9,
    7049.7,
    7097.6,
    7005.7,
    7045.4,
    7034.1,
    7067.3,
    6932.8,
    6888.8,
    6865.1,
    6828,
    6754.4,
    6738.5,
    6740.1,
    6747.8,
    6380.9,
    6405.6,
    6424.3,
    6465.4,
    6483,
    6538,
    6574.6,
    6639.2,
    6611.6,
    6613.9,
    6676.1,
    6738.9,
    6851.9,
    6853.3,
    6855.1,
    6882.8,
    7065,
    7020.6,
    7014.4,
    7005.6,
    7503.1,
    7551.6,
    7606,
    7655.5,
    7700,
    7750.5,
    7813,
    7833.9,
    8179.8,
    8181.1,
    8245.700000000001,
    8318.6,
    8433.6,
    8540.299999999999,
    8616,
    8684.9,
    8443.9,
    8468.1,
    8500.299999999999,
    8555.299999999999,
    8414.799999999999,
    8346.700000000001,
    8335.200000000001,
    8334.1,
    8325.5,
    8301.700000000001,
    8201,
    8208.700000000001,
    8169.8,
    8076.7,
    8058.9,
    8035,
    6282,
    6242,
    6182,
    6122,
    6126,
    6101,
    6055.8,
    6016.7,
    5995.2,
    5962.5,
    5977.3,
    5943.7,
    6014.3,
    6017.8,
    6022.4,
    6027.9,
    5996.8,
    6003.9,
    5994.1,
    6007.9,
    5936.6,
    5959.9,
    5994.8,
    6002.9,
    6017.8,
    6101.2,
    6198.7,
    6251,
    6928.9,
    6960.9,
    7069.7,
    7097.6,
    6995.7,
    7025.4,
    7034.1,
    7057.3,
    6892.8,
    6828.8,
    6775.1,
    6708,
    6594.4,
    6528.5,
    6490.1,
    6467.8,
    6060.9,
    6065.6,
    6064.3,
    6095.4,
    6103,
    6148,
    6174.6,
    6229.2,
    6191.6,
    6183.9,
    6236.1,
    6288.9,
    6411.9,
    6433.3,
    6465.1,
    6522.8,
    6735,
    6740.6,
    6764.4,
    6775.6,
    7293.1,
    7361.6,
    7426,
    7495.5,
    7560,
    7650.5,
    7753,
    7833.9,
    8219.799999999999,
    8261.1,
    8355.700000000001,
    8448.6,
    8563.6,
    8640.299999999999,
    8716,
    8744.9,
    8493.9,
    8558.1,
    8600.299999999999,
    8645.299999999999,
    8504.799999999999,
    8456.700000000001,
    8435.200000000001,
    8444.1,
    8425.5,
    8401.700000000001,
    8321,
    8308.700000000001,
    8249.799999999999,
    8136.7,
    8138.9,
    8115
)

$startRow = 2018
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newValues[$i]
}
